$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added address column in excel sheets"
# Insert a new blank column at F; this pushes the existing "District" column
# (previously F) one place to the right, to G, matching the diff.
$ws.Columns("F:F").Insert()

# New column header
$ws.Cells.Item(2, 6).Value = "Address"

# Populate the new "Address" column (F) for every data row. A handful of rows
# (3, 17, 28, 56, and the wrapped continuation rows 29/57 already carry their
# own data in B/E) are intentionally left blank, matching the source diff.
$ws.Cells.Item(4, 6).Value = 'S F D G H S HattimatturSavnur'
$ws.Cells.Item(5, 6).Value = 'B R E T English Medium High SchoolMotebennurByadgi'
$ws.Cells.Item(6, 6).Value = 'S V R P G H S KadashettihalliHangal'
$ws.Cells.Item(7, 6).Value = 'NTP GHS BudapanahalliByadgi'
$ws.Cells.Item(8, 6).Value = 'G H P S No – 8 Nagendranamatti'
$ws.Cells.Item(9, 6).Value = 'Sri Halasiddeshwar High SchoolHalageriRanebennur'
$ws.Cells.Item(10, 6).Value = 'Govt. H P S Haleritti'
$ws.Cells.Item(11, 6).Value = 'Sri Basaveshwara High SchoolKodiyalRanebenur'
$ws.Cells.Item(12, 6).Value = 'Govt. High School MantrodiSavanur'
$ws.Cells.Item(13, 6).Value = 'Shri Maruthi High School Hirekerur'
$ws.Cells.Item(14, 6).Value = 'G H S ChatraByadgi'
$ws.Cells.Item(15, 6).Value = 'G H S HalemanaggiSavanur'
$ws.Cells.Item(16, 6).Value = 'Govt. High SchoolGangapuraRanebennuru'
$ws.Cells.Item(18, 6).Value = 'S K H S KotihalRanebennur'
$ws.Cells.Item(19, 6).Value = 'G H S HirebidariRanebennur'
$ws.Cells.Item(20, 6).Value = 'NCJC High SchoolHangal'
$ws.Cells.Item(21, 6).Value = 'G H S NoolageriHirekerur'
$ws.Cells.Item(22, 6).Value = 'G H S Basapur'
$ws.Cells.Item(23, 6).Value = 'Hirekerur'
$ws.Cells.Item(24, 6).Value = 'SJJM Govt PU College (H S) Byadgi'
$ws.Cells.Item(25, 6).Value = 'G U H S MasurHirekerur'
$ws.Cells.Item(26, 6).Value = 'S S H S BalambeedHangal'
$ws.Cells.Item(27, 6).Value = 'Govt. High School HotanahalliShiggoan'
$ws.Cells.Item(29, 6).Value = 'G H S ItagiRanebennur'
$ws.Cells.Item(30, 6).Value = 'M D R S MaranabeeduHangal'
$ws.Cells.Item(31, 6).Value = 'S J J M Govt. PU CollegeByadgi'
$ws.Cells.Item(32, 6).Value = 'G H S KopparsikoppaHanagal'
$ws.Cells.Item(33, 6).Value = 'G H S HoovinashigliSavanur'
$ws.Cells.Item(34, 6).Value = 'G H P S Hounsi'
$ws.Cells.Item(35, 6).Value = 'S G B D H S ShidenurBayadagi'
$ws.Cells.Item(36, 6).Value = 'Govt Urdu High SchoolChikkerurHirekerur'
$ws.Cells.Item(37, 6).Value = 'G H P S RamagondanahalliByadgi'
$ws.Cells.Item(38, 6).Value = 'Govt. High School ShirabadagiSavanur'
$ws.Cells.Item(39, 6).Value = 'Sri Manjunath H S'
$ws.Cells.Item(40, 6).Value = 'Govt. High SchoolKalasurSavanur'
$ws.Cells.Item(41, 6).Value = 'Govt. High School kyalkondaShiggaon'
$ws.Cells.Item(42, 6).Value = 'P V S High SchoolChikkerurHirekerur'
$ws.Cells.Item(43, 6).Value = 'Vidyabhrati High School Savanur'
$ws.Cells.Item(44, 6).Value = 'Shri Kantesh High School KadaramandalagiByadagi'
$ws.Cells.Item(45, 6).Value = 'G H S BasavanalaShiggaon'
$ws.Cells.Item(46, 6).Value = 'S M H SchoolGhalapujiByadagi'
$ws.Cells.Item(47, 6).Value = 'Govt. Higher Primary School Naganur'
$ws.Cells.Item(48, 6).Value = 'Govt. High SchoolMakanurRanebennur'
$ws.Cells.Item(49, 6).Value = 'G H S NesviHirekerur'
$ws.Cells.Item(50, 6).Value = 'G H P S HalemannangiSavanur'
$ws.Cells.Item(51, 6).Value = 'G H S Kabbur'
$ws.Cells.Item(52, 6).Value = 'Sri Maruti High School KodHirekerur'
$ws.Cells.Item(53, 6).Value = 'Raj – RajeshwariHigh SchoolRanebennur'
$ws.Cells.Item(54, 6).Value = 'Govt. High SchoolKummurByadgi'
$ws.Cells.Item(55, 6).Value = 'G H S ShirgodHangal'
$ws.Cells.Item(57, 6).Value = 'G H S HullattiHirekerur'
$ws.Cells.Item(58, 6).Value = 'G H S HosahalliSavanur'
$ws.Cells.Item(59, 6).Value = 'S A H S GhalapujiBydagi'
$ws.Cells.Item(60, 6).Value = 'S S M H SchoolGuddadamallapurByadgi'
$ws.Cells.Item(61, 6).Value = 'S G R High SchoolSavanur'
$ws.Cells.Item(62, 6).Value = 'G H S HireanajiByadgi'
$ws.Cells.Item(63, 6).Value = 'Govt. High SchoolKeravadiByadgi'
$ws.Cells.Item(64, 6).Value = 'Govt. High SchoolMasanagiByadagi'
$ws.Cells.Item(65, 6).Value = 'G H S HunagundShiggaon'
$ws.Cells.Item(66, 6).Value = 'S S P U College TumminakattiRanebennur'
$ws.Cells.Item(67, 6).Value = 'G H P S Guddada ChannapurShiggaon'
$ws.Cells.Item(68, 6).Value = 'SJJM Govt PU CollegeByadgi'
$ws.Cells.Item(69, 6).Value = 'SSSPUC (HS) SunakalbidriRanebennuru'

